$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 385
$wsExpo.Range("F4").Value = 4920
$wsExpo.Range("F5").Value = 26

# Sheet "全部类型" (All types) - same underlying events, update matching cells
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 385
$wsAll.Range("F4").Value = 4920
$wsAll.Range("F6").Value = 26
